# Regenerate orders with updated distance/size codes.
# The experiment's "distance" and "size" condition labels changed:
#   D80 -> D86
#   D51 -> D55
#   D64 -> D69
#   S30 -> S31
# These codes appear embedded throughout several columns (Condition,
# Filename_Left, Filename_Right, Distance, Size), so perform a global
# find/replace across the sheet's used range for each code, matching
# whole cell text occurrences (Range.Replace mirrors Excel's Ctrl+H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.UsedRange

$rng.Replace("D80", "D86")
$rng.Replace("D51", "D55")
$rng.Replace("D64", "D69")
$rng.Replace("S30", "S31")
